# New crime data collected — weekly CompStat refresh (112th Precinct):
#   - bump "Volume 31   Number 51" -> "...Number 52" in the banner
#   - advance the reporting week 12/16/2024-12/22/2024 -> 12/23/2024-12/29/2024
#   - widen column I to match column H's width (bestFit growth)
#   - refresh the crime-stat grid (rows 15-31) with the new week's numbers,
#     including a few cells flipping between the literal "0"/"***.*" placeholder
#     text and real numbers as categories move in/out of having data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header banner text — edit just the trailing digits/date via Characters()
# so the rest of the rich-text run (font/size/color) is left alone.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "52"

$ws.Range("C9").Characters(27, 10).Text = "12/23/2024"
$ws.Range("C9").Characters(48, 10).Text = "12/29/2024"

# ---------------------------------------------------------------------------
# Column I grows to match column H's bestFit width (7.433768 chars). The
# host's ColumnWidth setter quantises to 1/7-character steps, so 6.71 is the
# closest achievable value (-> stored width ~7.4286, H's own read-back value).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 6.71

# ---------------------------------------------------------------------------
# Plain numeric updates — style/type unchanged, so a direct .Value assign
# is enough.
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 133.333333333333
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 40

$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 87
$ws.Range("K16").Value = 20.833333333333
$ws.Range("L16").Value = -1.136363636363
$ws.Range("M16").Value = -26.890756302521
$ws.Range("N16").Value = -86.778115501519

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 114
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = 5.555555555555
$ws.Range("L17").Value = 34.117647058823
$ws.Range("M17").Value = 96.551724137931
$ws.Range("N17").Value = -8.8

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 13
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -9.243697478991
$ws.Range("L18").Value = -10.743801652892
$ws.Range("M18").Value = -20.588235294117
$ws.Range("N18").Value = -92.207792207792

$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -11.428571428571
$ws.Range("I19").Value = 449
$ws.Range("J19").Value = 479
$ws.Range("K19").Value = -6.263048016701
$ws.Range("L19").Value = -9.109311740890
$ws.Range("M19").Value = 16.623376623376
$ws.Range("N19").Value = -56.023506366307

$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -58.333333333333
$ws.Range("I20").Value = 161
$ws.Range("J20").Value = 170
$ws.Range("K20").Value = -5.294117647058
$ws.Range("L20").Value = 51.886792452830
$ws.Range("M20").Value = 62.626262626262
$ws.Range("N20").Value = -95.221133867616

$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -7.594936708860
$ws.Range("I21").Value = 933
$ws.Range("J21").Value = 954
$ws.Range("K21").Value = -2.201257861635
$ws.Range("L21").Value = 2.302631578947
$ws.Range("M21").Value = 16.334164588528
$ws.Range("N21").Value = -85.807727411013

$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 39
$ws.Range("K22").Value = 39.285714285714
$ws.Range("L22").Value = 8.333333333333
$ws.Range("M22").Value = 77.272727272727

$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = 5.405405405405
$ws.Range("F24").Value = 160
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = 19.402985074626
$ws.Range("I24").Value = 1788
$ws.Range("J24").Value = 1553
$ws.Range("K24").Value = 15.13200257566
$ws.Range("L24").Value = 0.393037619314
$ws.Range("M24").Value = 80.788675429727

$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -3.571428571428
$ws.Range("F25").Value = 115
$ws.Range("G25").Value = 96
$ws.Range("H25").Value = 19.791666666666
$ws.Range("I25").Value = 1317
$ws.Range("J25").Value = 1108
$ws.Range("K25").Value = 18.862815884476
$ws.Range("L25").Value = 4.440919904837

$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -20.833333333333
$ws.Range("I26").Value = 274
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 9.6
$ws.Range("L26").Value = 41.237113402061
$ws.Range("M26").Value = 27.441860465116

$ws.Range("F27").Value = 4
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = 76.923076923076
$ws.Range("L27").Value = 4.545454545454

$ws.Range("J31").Value = 15
$ws.Range("K31").Value = -26.666666666666

# ---------------------------------------------------------------------------
# Cells flipping from a number to the shared "0" / "***.*" placeholder text.
# A bare .Value assign of a numeric-looking string ("0") gets auto-coerced
# back to a number, and any string assign also forces a fresh (non "s=13")
# style — so nail the value down first, then paste just the formatting from
# a same-shaped neighbour (row 23, untouched by this week's edit) to restore
# the shared "text" style (s=13) these placeholder cells use elsewhere.
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("D16").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("G27").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("G27").PasteSpecial(-4122)

$ws.Range("H27").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("H27").PasteSpecial(-4122)

$ws.Range("C28").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Cells flipping the other way: placeholder text -> a real number. Paste the
# number/percent formatting back from row 23's equivalent column so the
# style returns to the plain numeric style (s=15) / percent style (s=14)
# used throughout the rest of the grid.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("I23").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D31").Value = 1
$ws.Range("I23").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("G31").Value = 1
$ws.Range("I23").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("E31").Value = -100
$ws.Range("K23").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("H31").Value = -100
$ws.Range("K23").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = $false
